# AFDP-2369: update the "next possible queue" rule in the Drools decision
# table on Sheet1 and correct a stray indentation typo in the helper
# function text used by the rule.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The CONDITION column previously matched the wildcard queue "*"; replace it
# with a condition that checks for a specific (non-matching) queue name.
$ws.Range("C18").Value = 'queue.name.equals("no-such-queue")'

# Fix the mis-indented line inside the evalSpring() helper function body
# (single leading space -> four spaces, matching the rest of the function).
$evalSpringFunction = "function Boolean evalSpring(String expression, NextPossibleQueuesModel model)`n{`n    ExpressionParser ep = new SpelExpressionParser();`n    Expression exp = ep.parseExpression(expression);`n    EvaluationContext ec = new StandardEvaluationContext();`n `n    CaseFile caseFile = (CaseFile) model.getBusinessObject();`n`n    Boolean evaluated = exp.getValue(ec, caseFile, Boolean.class);`n `n    return evaluated;`n}"
$ws.Range("D10").Value = $evalSpringFunction

# Update the sheet's view/selection state to match the edited cell.
[void]$ws.Activate()
[void]$ws.Range("C18").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1

Write-Host "C18 is now: $($ws.Range("C18").Value())"
